# daily auto push: 2026-01-07 18:49 UTC
# Two new sensor readings land for 2026/01/07-08, which need to be inserted
# in their correct chronological position (between the existing 2026/01/07
# rows and the 2026/12/29 rows), pushing the rest of the table down by two
# rows. Two more readings for 2027/01/05 are appended at the very end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows right before the old row 588 (the first row of the
# 2026/12/29 block) - this shifts the existing 588:629 block down to 590:631
# and leaves 588:589 empty for the new readings.
$ws.Rows("588:589").Insert()

# Column A holds dates formatted as plain text (e.g. "2026/01/07"); mark the
# new cells as text first so Excel does not silently convert them into date
# serial numbers, then restore the default "Normal" style so the cell ends
# up styled exactly like its neighbours (no leftover text-format override).
$ws.Range("A588:A589").NumberFormat = "@"
$ws.Range("A630:A631").NumberFormat = "@"

# New row 588: 2026/01/07 23:00
$ws.Cells.Item(588, 1).Value = "2026/01/07"
$ws.Cells.Item(588, 2).Value = "水"
$ws.Cells.Item(588, 3).Value = 23
$ws.Cells.Item(588, 4).Value = 201

# New row 589: 2026/01/08 02:00
$ws.Cells.Item(589, 1).Value = "2026/01/08"
$ws.Cells.Item(589, 2).Value = "木"
$ws.Cells.Item(589, 3).Value = 2
$ws.Cells.Item(589, 4).Value = 201

# Appended row 630: 2027/01/05 02:00
$ws.Cells.Item(630, 1).Value = "2027/01/05"
$ws.Cells.Item(630, 2).Value = "火"
$ws.Cells.Item(630, 3).Value = 2
$ws.Cells.Item(630, 4).Value = 201

# Appended row 631: 2027/01/05 07:00
$ws.Cells.Item(631, 1).Value = "2027/01/05"
$ws.Cells.Item(631, 2).Value = "火"
$ws.Cells.Item(631, 3).Value = 7
$ws.Cells.Item(631, 4).Value = 201

# Drop the temporary text-format override now that the literal strings are
# safely stored, so these cells match the unstyled look of the rest of the
# table.
$ws.Range("A588:A589").Style = "Normal"
$ws.Range("A630:A631").Style = "Normal"
